# edit.ps1
# Applies updated category usage-frequency data (rows 2-52) to Sheet1,
# matching the refreshed/resorted filtered data described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('outros', 154, 590, 0, 0),
    @('saúde', 24, 61, 4, 0),
    @('sistema penitenciário', 3, 22, 2, 1),
    @('nacionalidade e nacionalismo', 27, 66, 0, 1),
    @('vida rural, vida no campo', 21, 79, 5, 1),
    @('amizade', 4, 23, 3, 3),
    @('pessoas com deficiência', 8, 9, 3, 3),
    @('maternidade', 5, 15, 2, 4),
    @('raça e racismo', 3, 23, 6, 4),
    @('sociedade e consumo', 59, 202, 6, 5),
    @('linguagem e representação', 81, 240, 5, 5),
    @('violência', 31, 96, 8, 5),
    @('recreação, lazer e entretenimento', 138, 462, 9, 6),
    @('família', 27, 162, 44, 6),
    @('sonho e fantasia', 11, 22, 2, 6),
    @('desinformação, populismo e polarização', 39, 134, 7, 7),
    @('distopia, ficção científica e futuros imaginados', 68, 121, 5, 7),
    @('tecnologia, inovação e sociedade', 33, 114, 13, 7),
    @('questão agrária e territorial', 33, 127, 16, 8),
    @('moradia e habitação', 6, 18, 21, 8),
    @('mobilidade', 31, 102, 21, 10),
    @('territorialidade e colonialismo', 34, 86, 16, 11),
    @('violências e preconceitos de gênero', 16, 70, 4, 12),
    @('alimentação e tratamentos tradicionais', 31, 116, 8, 13),
    @('sonoridade e paisagem sonora', 49, 140, 11, 14),
    @('gênero e sexualidade', 70, 322, 10, 14),
    @('biografia', 10, 25, 0, 15),
    @('trabalho e ofício', 86, 286, 29, 16),
    @('direitos humanos', 26, 88, 9, 17),
    @('crises e desastres ambientais e sociais', 29, 86, 14, 17),
    @('religião, espiritualidade e cosmologias', 38, 124, 14, 19),
    @('solidão', 9, 27, 15, 19),
    @('saudade, luto ou perda', 51, 110, 7, 20),
    @('política', 86, 380, 29, 23),
    @('dinâmica urbana', 168, 427, 48, 23),
    @('vida afetiva', 313, 1001, 59, 25),
    @('estruturas sociais e econômicas', 159, 520, 23, 26),
    @('educação e socialização', 104, 312, 26, 26),
    @('cultural', 196, 660, 43, 27),
    @('corpo, performance e expressão', 33, 156, 6, 27),
    @('memória e patrimônio', 114, 379, 25, 27),
    @('povos originários e comunidades tradicionais', 76, 228, 19, 29),
    @('pandemia', 113, 304, 48, 29),
    @('mulher, feminino e feminismo', 21, 61, 22, 31),
    @('vida cotidiana', 24, 107, 18, 37),
    @('resistência e luta', 22, 80, 27, 37),
    @('saúde mental', 106, 284, 37, 39),
    @('ambiental', 114, 346, 52, 46),
    @('reflexão', 9, 41, 19, 47),
    @('poesia e ensaio', 223, 696, 56, 50),
    @('arte', 20, 41, 24, 61)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

$wb.Save()
